# Append: 2025-11-08 01:43 JST
# Update the "取得日時" (fetched-at) timestamp column (A2:A15) on the
# first worksheet from the old scrape time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "2025-11-08 01:15:48"
$newValue = "2025-11-08 01:43:15"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
